$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as a literal text string into a cell without
# Excel's automatic "looks like a number" reinterpretation (which would
# silently strip things like trailing zeros, e.g. "217.90" -> 217.9).
# We do this by writing a formula that evaluates to the literal text,
# then converting the whole range to static values via copy/paste-special
# (values only) so the final cells end up as plain text values with no
# formulas and no style changes.
function Set-TextValue($row, $col, $text) {
    $escaped = $text.Replace('"', '""')
    $ws.Cells.Item($row, $col).Formula = '="' + $escaped + '"'
}

Set-TextValue 2 4 "91.232.13"
Set-TextValue 2 5 "  +3.63%  "
Set-TextValue 3 4 "3.096.67"
Set-TextValue 3 5 "  -0.39%  "
Set-TextValue 4 5 "  +0.00%  "
Set-TextValue 5 4 "217.90"
Set-TextValue 5 5 "  +2.14%  "
Set-TextValue 6 4 "618.17"
Set-TextValue 6 5 "  -2.36%  "
Set-TextValue 7 4 "0.377"
Set-TextValue 7 5 "  -0.38%  "
Set-TextValue 8 4 "0.904"
Set-TextValue 8 5 "  +11.55%  "
Set-TextValue 9 5 "  -0.04%  "
Set-TextValue 10 4 "3.092.35"
Set-TextValue 10 5 "  -0.41%  "
Set-TextValue 11 4 "0.672"
Set-TextValue 11 5 "  +14.54%  "
Set-TextValue 12 5 "  +5.99%  "
Set-TextValue 13 4 "0.0000254"
Set-TextValue 13 5 "  +4.02%  "
Set-TextValue 14 4 "5.39"
Set-TextValue 14 5 "  +0.86%  "
Set-TextValue 15 4 "90.945.90"
Set-TextValue 15 5 "  +3.58%  "
Set-TextValue 16 5 "  +3.66%  "
Set-TextValue 17 4 "3.673.78"
Set-TextValue 17 5 "  +0.00%  "
Set-TextValue 18 4 "3.059.41"
Set-TextValue 18 5 "  -1.35%  "
Set-TextValue 19 4 "3.66"
Set-TextValue 19 5 "  +6.49%  "
Set-TextValue 20 4 "0.0000223"
Set-TextValue 20 5 "  +4.40%  "
Set-TextValue 21 4 "13.84"
Set-TextValue 21 5 "  +4.14%  "
Set-TextValue 22 4 "432.30"
Set-TextValue 22 5 "  +2.20%  "
Set-TextValue 23 4 "8.53"
Set-TextValue 23 5 "  +1.62%  "
Set-TextValue 24 4 "5.13"
Set-TextValue 24 5 "  +5.77%  "
Set-TextValue 25 5 "  +1.56%  "
Set-TextValue 26 2 "Aptos"
Set-TextValue 26 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue 26 4 "11.89"
Set-TextValue 26 5 "  +4.15%  "
Set-TextValue 27 2 "Litecoin"
Set-TextValue 27 3 "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue 27 4 "83.99"
Set-TextValue 27 5 "  +0.24%  "
Set-TextValue 28 4 "3.261.68"
Set-TextValue 28 5 "  -0.25%  "
Set-TextValue 29 4 "0.998"
Set-TextValue 29 5 "  -0.18%  "
Set-TextValue 30 5 "  +10.33%  "
Set-TextValue 31 4 "1.00"
Set-TextValue 31 5 "  +0.00%  "
Set-TextValue 32 5 "  +7.23%  "
Set-TextValue 33 4 "3.87"
Set-TextValue 33 5 "  +2.49%  "
Set-TextValue 34 4 "519.14"
Set-TextValue 34 5 "  +2.96%  "
Set-TextValue 35 4 "7.04"
Set-TextValue 35 5 "  +3.50%  "
Set-TextValue 36 5 "  -2.46%  "
Set-TextValue 37 5 "  +1.89%  "
Set-TextValue 38 5 "  +1.16%  "
Set-TextValue 39 4 "23.02"
Set-TextValue 39 5 "  +2.53%  "
Set-TextValue 40 4 "22.29"
Set-TextValue 40 5 "  +0.51%  "
Set-TextValue 41 5 "  -0.10%  "
Set-TextValue 43 4 "0.143"
Set-TextValue 43 5 "  +5.30%  "
Set-TextValue 44 4 "0.370"
Set-TextValue 44 5 "  +1.18%  "
Set-TextValue 45 4 "1.87"
Set-TextValue 45 5 "  +1.48%  "
Set-TextValue 46 4 "0.0727"
Set-TextValue 46 5 "  +9.98%  "
Set-TextValue 47 2 "Monero"
Set-TextValue 47 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 47 4 "142.56"
Set-TextValue 47 5 "  -2.47%  "
Set-TextValue 48 2 "OKB"
Set-TextValue 48 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue 48 4 "43.55"
Set-TextValue 48 5 "  -0.17%  "
Set-TextValue 49 4 "0.000260"
Set-TextValue 49 5 "  +14.34%  "
Set-TextValue 50 4 "4.19"
Set-TextValue 50 5 "  +6.09%  "
Set-TextValue 51 4 "1.24"
Set-TextValue 51 5 "  +4.60%  "

# Convert the formulas above into plain static text values, preserving
# the original cell styles/number formats (no new styles are introduced).
$dataRange = $ws.Range("A1:E51")
$dataRange.Copy()
$dataRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false
